$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("region") held bare integer region indices (1..78) in rows 3:80.
# Replace each with a descriptive label "deg_region_<n>" so the region
# identifier is self-explanatory without needing the header for context.
for ($row = 3; $row -le 80; $row++) {
    $regionNumber = $row - 2
    $ws.Cells.Item($row, 2).Value = "deg_region_$regionNumber"
}
